$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (in-place run edits on the shared/rich-text strings) ---
# A8 holds "Volume 30   Number  6" as 4 runs; only the trailing "6" run changes to "7".
$ws.Range("A8").Characters(21, 1).Text = "7"
# C9 holds "Report Covering the Week  2/6/2023  Through  2/12/2023" as 4 runs;
# the two date runs change (lengths differ, so do the earlier edit first).
$ws.Range("C9").Characters(27, 8).Text = "2/13/2023"
$ws.Range("C9").Characters(47, 9).Text = "2/19/2023"

# --- Weekly crime statistics table updates (rows 14-30) ---
# Row 14
$ws.Range("C14").Value = 2
$ws.Range("D14").Value = 2
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 10
$ws.Range("G14").Value = 10
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 15
$ws.Range("J14").Value = 15
$ws.Range("L14").Value = 7.142857142857
$ws.Range("M14").Value = 87.5
$ws.Range("N14").Value = -76.5625

# Row 15
$ws.Range("C15").Value = 8
$ws.Range("D15").Value = 7
$ws.Range("E15").Value = 14.285714285714
$ws.Range("F15").Value = 28
$ws.Range("G15").Value = 32
$ws.Range("H15").Value = -12.5
$ws.Range("I15").Value = 54
$ws.Range("J15").Value = 61
$ws.Range("K15").Value = -11.475409836065
$ws.Range("L15").Value = 31.70731707317
$ws.Range("M15").Value = 45.945945945945
$ws.Range("N15").Value = -22.857142857142

# Row 16
$ws.Range("C16").Value = 66
$ws.Range("D16").Value = 84
$ws.Range("E16").Value = -21.428571428571
$ws.Range("F16").Value = 294
$ws.Range("G16").Value = 299
$ws.Range("H16").Value = -1.672240802675
$ws.Range("I16").Value = 582
$ws.Range("J16").Value = 572
$ws.Range("K16").Value = 1.748251748251
$ws.Range("L16").Value = 31.081081081081
$ws.Range("M16").Value = 5.818181818181
$ws.Range("N16").Value = -74.585152838427

# Row 17
$ws.Range("C17").Value = 134
$ws.Range("D17").Value = 114
$ws.Range("E17").Value = 17.543859649122
$ws.Range("F17").Value = 505
$ws.Range("H17").Value = 9.307359307359
$ws.Range("I17").Value = 912
$ws.Range("J17").Value = 826
$ws.Range("K17").Value = 10.411622276029
$ws.Range("L17").Value = 27.910238429172
$ws.Range("M17").Value = 71.751412429378
$ws.Range("N17").Value = -6.843718079673

# Row 18
$ws.Range("C18").Value = 65
$ws.Range("D18").Value = 58
$ws.Range("E18").Value = 12.068965517241
$ws.Range("F18").Value = 242
$ws.Range("G18").Value = 221
$ws.Range("H18").Value = 9.502262443438
$ws.Range("I18").Value = 417
$ws.Range("J18").Value = 383
$ws.Range("K18").Value = 8.8772845953
$ws.Range("L18").Value = 57.358490566037
$ws.Range("M18").Value = -6.919642857142
$ws.Range("N18").Value = -83.537307540465

# Row 19
$ws.Range("C19").Value = 115
$ws.Range("D19").Value = 170
$ws.Range("E19").Value = -32.35294117647
$ws.Range("F19").Value = 505
$ws.Range("G19").Value = 578
$ws.Range("H19").Value = -12.629757785467
$ws.Range("I19").Value = 936
$ws.Range("J19").Value = 1041
$ws.Range("K19").Value = -10.086455331412
$ws.Range("L19").Value = 31.460674157303
$ws.Range("M19").Value = 78.967495219885
$ws.Range("N19").Value = 4

# Row 20
$ws.Range("C20").Value = 93
$ws.Range("D20").Value = 84
$ws.Range("E20").Value = 10.714285714285
$ws.Range("F20").Value = 405
$ws.Range("G20").Value = 347
$ws.Range("H20").Value = 16.71469740634
$ws.Range("I20").Value = 740
$ws.Range("J20").Value = 643
$ws.Range("K20").Value = 15.085536547433
$ws.Range("L20").Value = 183.524904214559
$ws.Range("M20").Value = 184.615384615385
$ws.Range("N20").Value = -65.061378659112

# Row 21
$ws.Range("C21").Value = 483
$ws.Range("D21").Value = 519
$ws.Range("E21").Value = -6.936416184971
$ws.Range("F21").Value = 1989
$ws.Range("G21").Value = 1949
$ws.Range("H21").Value = 2.052334530528
$ws.Range("I21").Value = 3656
$ws.Range("J21").Value = 3541
$ws.Range("K21").Value = 3.247670149675
$ws.Range("L21").Value = 49.224489795918
$ws.Range("M21").Value = 55.112431056427
$ws.Range("N21").Value = -59.169086441813

# Row 22
$ws.Range("C22").Value = 9
$ws.Range("D22").Value = 6
$ws.Range("E22").Value = 50
$ws.Range("F22").Value = 21
$ws.Range("G22").Value = 26
$ws.Range("H22").Value = -19.230769230769
$ws.Range("I22").Value = 32
$ws.Range("J22").Value = 48
$ws.Range("K22").Value = -33.333333333333
$ws.Range("L22").Value = 14.285714285714
$ws.Range("M22").Value = 0

# Row 23
$ws.Range("C23").Value = 28
$ws.Range("D23").Value = 31
$ws.Range("E23").Value = -9.677419354838
$ws.Range("F23").Value = 123
$ws.Range("G23").Value = 121
$ws.Range("H23").Value = 1.652892561983
$ws.Range("I23").Value = 230
$ws.Range("J23").Value = 200
$ws.Range("K23").Value = 15
$ws.Range("L23").Value = 57.534246575342
$ws.Range("M23").Value = 75.572519083969

# Row 24
$ws.Range("C24").Value = 358
$ws.Range("D24").Value = 325
$ws.Range("E24").Value = 10.153846153846
$ws.Range("F24").Value = 1285
$ws.Range("G24").Value = 1261
$ws.Range("H24").Value = 1.903251387787
$ws.Range("I24").Value = 2213
$ws.Range("J24").Value = 2113
$ws.Range("K24").Value = 4.732607666824
$ws.Range("L24").Value = 36.94306930693
$ws.Range("M24").Value = 40.329740012682

# Row 25
$ws.Range("C25").Value = 183
$ws.Range("D25").Value = 162
$ws.Range("E25").Value = 12.962962962963
$ws.Range("F25").Value = 752
$ws.Range("G25").Value = 724
$ws.Range("H25").Value = 3.867403314917
$ws.Range("I25").Value = 1276
$ws.Range("J25").Value = 1194
$ws.Range("K25").Value = 6.867671691792
$ws.Range("L25").Value = 31.546391752577
$ws.Range("M25").Value = 3.65556458164

# Row 26
$ws.Range("C26").Value = 12
$ws.Range("D26").Value = 11
$ws.Range("E26").Value = 9.090909090909
$ws.Range("F26").Value = 48
$ws.Range("G26").Value = 56
$ws.Range("H26").Value = -14.285714285714
$ws.Range("I26").Value = 89
$ws.Range("J26").Value = 98
$ws.Range("K26").Value = -9.183673469387
$ws.Range("L26").Value = 34.848484848484

# Row 27
$ws.Range("C27").Value = 23
$ws.Range("D27").Value = 16
$ws.Range("E27").Value = 43.75
$ws.Range("F27").Value = 68
$ws.Range("G27").Value = 68
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 151
$ws.Range("J27").Value = 105
$ws.Range("K27").Value = 43.809523809523
$ws.Range("L27").Value = 36.036036036036

# Row 28
$ws.Range("C28").Value = 3
$ws.Range("D28").Value = 4
$ws.Range("E28").Value = -25
$ws.Range("F28").Value = 22
$ws.Range("G28").Value = 39
$ws.Range("H28").Value = -43.589743589743
$ws.Range("I28").Value = 41
$ws.Range("J28").Value = 62
$ws.Range("K28").Value = -33.870967741935
$ws.Range("L28").Value = -2.380952380952
$ws.Range("M28").Value = -18
$ws.Range("N28").Value = -71.917808219178

# Row 29
$ws.Range("C29").Value = 3
$ws.Range("D29").Value = 4
$ws.Range("E29").Value = -25
$ws.Range("F29").Value = 17
$ws.Range("G29").Value = 34
$ws.Range("H29").Value = -50
$ws.Range("I29").Value = 32
$ws.Range("J29").Value = 57
$ws.Range("K29").Value = -43.859649122807
$ws.Range("L29").Value = -15.78947368421
$ws.Range("M29").Value = -28.888888888888
$ws.Range("N29").Value = -75.757575757575

# Row 30
$ws.Range("D30").Value = 3
$ws.Range("G30").Value = 7
$ws.Range("H30").Value = -100
$ws.Range("J30").Value = 10
$ws.Range("K30").Value = -60

# F30 switches from a numeric 1 to the literal text "0" (same style used by C30).
# Paste C30s format onto F30 (keeps it on style s="14"), then paste its value
# (the shared string "0") without disturbing the format just applied.
$ws.Range("C30").Copy()
$ws.Range("F30").PasteSpecial(-4122)
$ws.Range("C30").Copy()
$ws.Range("F30").PasteSpecial(-4163)
$excel.CutCopyMode = $false
